$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 224.6842
$ws.Cells.Item(33, 9).Value = 282.15384
$ws.Cells.Item(33, 10).Value = 100.166664
$ws.Cells.Item(33, 11).Value = 282.15384
$ws.Cells.Item(33, 12).Value = 100.166664
$ws.Cells.Item(33, 13).Value = -53.15384
$ws.Cells.Item(33, 14).Value = -558.166664
$ws.Cells.Item(51, 8).Value = 2868.8462
$ws.Cells.Item(51, 10).Value = 3077.8572
$ws.Cells.Item(51, 12).Value = 3077.8572
$ws.Cells.Item(51, 14).Value = -4045.8572
$ws.Cells.Item(52, 8).Value = 1283.7778
$ws.Cells.Item(52, 9).Value = 1264.5
$ws.Cells.Item(52, 10).Value = 1299.2
$ws.Cells.Item(52, 11).Value = 3793.5
$ws.Cells.Item(52, 12).Value = 3897.6
$ws.Cells.Item(52, 13).Value = -3633.5
$ws.Cells.Item(52, 14).Value = -4217.6
$ws.Cells.Item(116, 8).Value = 7182.933
$ws.Cells.Item(116, 9).Value = 6580.4287
$ws.Cells.Item(116, 11).Value = 6580.4287
$ws.Cells.Item(116, 13).Value = -3138.4287
$ws.Cells.Item(125, 8).Value = 1500.0526
$ws.Cells.Item(125, 9).Value = 1059.7778
$ws.Cells.Item(125, 10).Value = 1896.3
$ws.Cells.Item(125, 11).Value = 9538.0002
$ws.Cells.Item(125, 12).Value = 17066.7
$ws.Cells.Item(125, 13).Value = -7078.0002
$ws.Cells.Item(125, 14).Value = -21986.7
$ws.Cells.Item(137, 8).Value = 1807.8649
$ws.Cells.Item(137, 9).Value = 1928.1428
$ws.Cells.Item(137, 10).Value = 1650
$ws.Cells.Item(137, 11).Value = 5784.428400000001
$ws.Cells.Item(137, 12).Value = 4950
$ws.Cells.Item(137, 13).Value = -3234.428400000001
$ws.Cells.Item(137, 14).Value = -10050
$ws.Cells.Item(138, 8).Value = 6759825
$ws.Cells.Item(138, 10).Value = 7939828.5
$ws.Cells.Item(138, 12).Value = 23819485.5
$ws.Cells.Item(138, 14).Value = -23829765.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7407.4165
$ws.Cells.Item(32, 9).Value = 3167.2913
$ws.Cells.Item(32, 11).Value = 3167.2913
$ws.Cells.Item(32, 13).Value = -2880.2913
$ws.Cells.Item(61, 8).Value = 5124.914
$ws.Cells.Item(61, 9).Value = 4114.3335
$ws.Cells.Item(61, 11).Value = 4114.3335
$ws.Cells.Item(61, 13).Value = -3902.3335
$ws.Cells.Item(74, 8).Value = 72088.94
$ws.Cells.Item(74, 9).Value = 126692.875
$ws.Cells.Item(74, 10).Value = 17485
$ws.Cells.Item(74, 11).Value = 126692.875
$ws.Cells.Item(74, 12).Value = 17485
$ws.Cells.Item(74, 13).Value = -125818.875
$ws.Cells.Item(74, 14).Value = -19233
$ws.Cells.Item(77, 8).Value = 72088.94
$ws.Cells.Item(77, 9).Value = 126692.875
$ws.Cells.Item(77, 10).Value = 17485
$ws.Cells.Item(77, 11).Value = 633464.375
$ws.Cells.Item(77, 12).Value = 87425
$ws.Cells.Item(77, 13).Value = -629096.375
$ws.Cells.Item(77, 14).Value = -96161
$ws.Cells.Item(102, 8).Value = 3800
$ws.Cells.Item(102, 9).Value = 3700
$ws.Cells.Item(102, 11).Value = 3700
$ws.Cells.Item(102, 13).Value = -2078
$ws.Cells.Item(122, 8).Value = 1166.9
$ws.Cells.Item(122, 9).Value = 608.75
$ws.Cells.Item(122, 11).Value = 1826.25
$ws.Cells.Item(122, 13).Value = 623.75
$ws.Cells.Item(132, 8).Value = 3460.239
$ws.Cells.Item(132, 9).Value = 3230.3157
$ws.Cells.Item(132, 10).Value = 4552.375
$ws.Cells.Item(132, 11).Value = 9690.947100000001
$ws.Cells.Item(132, 12).Value = 13657.125
$ws.Cells.Item(132, 13).Value = -7160.947100000001
$ws.Cells.Item(132, 14).Value = -18717.125
$ws.Cells.Item(136, 8).Value = 5124.914
$ws.Cells.Item(136, 9).Value = 4114.3335
$ws.Cells.Item(136, 11).Value = 12343.0005
$ws.Cells.Item(136, 13).Value = -9793.000499999998
$ws.Cells.Item(139, 8).Value = 63500
$ws.Cells.Item(139, 10).Value = 63500
$ws.Cells.Item(139, 12).Value = 63500
$ws.Cells.Item(139, 14).Value = -73780
$ws.Cells.Item(140, 8).Value = 100306.336
$ws.Cells.Item(140, 10).Value = 100306.336
$ws.Cells.Item(140, 12).Value = 100306.336
$ws.Cells.Item(140, 14).Value = -110666.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(58, 8).Value = 25907
$ws.Cells.Item(58, 9).Value = 9999
$ws.Cells.Item(58, 10).Value = 31209.666
$ws.Cells.Item(58, 11).Value = 9999
$ws.Cells.Item(58, 12).Value = 31209.666
$ws.Cells.Item(58, 13).Value = -9705
$ws.Cells.Item(58, 14).Value = -31797.666
$ws.Cells.Item(86, 8).Value = 9145.111000000001
$ws.Cells.Item(86, 10).Value = 7449.5
$ws.Cells.Item(86, 12).Value = 7449.5
$ws.Cells.Item(86, 14).Value = -9695.5
$ws.Cells.Item(89, 8).Value = 9145.111000000001
$ws.Cells.Item(89, 10).Value = 7449.5
$ws.Cells.Item(89, 12).Value = 37247.5
$ws.Cells.Item(89, 14).Value = -48479.5
$ws.Cells.Item(107, 8).Value = 1539.2273
$ws.Cells.Item(107, 9).Value = 1360.5264
$ws.Cells.Item(107, 10).Value = 2671
$ws.Cells.Item(107, 11).Value = 1360.5264
$ws.Cells.Item(107, 12).Value = 2671
$ws.Cells.Item(107, 13).Value = 559.4736
$ws.Cells.Item(107, 14).Value = -6511

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2824.1904
$ws.Cells.Item(58, 9).Value = 2361.625
$ws.Cells.Item(58, 10).Value = 3108.8462
$ws.Cells.Item(58, 11).Value = 2361.625
$ws.Cells.Item(58, 12).Value = 3108.8462
$ws.Cells.Item(58, 13).Value = -2158.625
$ws.Cells.Item(58, 14).Value = -3514.8462
$ws.Cells.Item(136, 8).Value = 2824.1904
$ws.Cells.Item(136, 9).Value = 2361.625
$ws.Cells.Item(136, 10).Value = 3108.8462
$ws.Cells.Item(136, 11).Value = 7084.875
$ws.Cells.Item(136, 12).Value = 9326.5386
$ws.Cells.Item(136, 13).Value = -4534.875
$ws.Cells.Item(136, 14).Value = -14426.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(29, 8).Value = 478.55
$ws.Cells.Item(29, 9).Value = 1088.375
$ws.Cells.Item(29, 10).Value = 72
$ws.Cells.Item(29, 11).Value = 1088.375
$ws.Cells.Item(29, 12).Value = 72
$ws.Cells.Item(29, 13).Value = -798.375
$ws.Cells.Item(29, 14).Value = -652
$ws.Cells.Item(70, 8).Value = 19505
$ws.Cells.Item(70, 9).Value = 14812.25
$ws.Cells.Item(70, 10).Value = 22633.5
$ws.Cells.Item(70, 11).Value = 14812.25
$ws.Cells.Item(70, 12).Value = 22633.5
$ws.Cells.Item(70, 13).Value = -14542.25
$ws.Cells.Item(70, 14).Value = -23173.5
$ws.Cells.Item(73, 8).Value = 19505
$ws.Cells.Item(73, 9).Value = 14812.25
$ws.Cells.Item(73, 10).Value = 22633.5
$ws.Cells.Item(73, 11).Value = 14812.25
$ws.Cells.Item(73, 12).Value = 22633.5
$ws.Cells.Item(73, 13).Value = -13876.25
$ws.Cells.Item(73, 14).Value = -24505.5
$ws.Cells.Item(80, 8).Value = 2998.75
$ws.Cells.Item(80, 9).Value = 2166.6667
$ws.Cells.Item(80, 10).Value = 3498
$ws.Cells.Item(80, 11).Value = 2166.6667
$ws.Cells.Item(80, 12).Value = 3498
$ws.Cells.Item(80, 13).Value = -1168.6667
$ws.Cells.Item(80, 14).Value = -5494
$ws.Cells.Item(83, 8).Value = 2998.75
$ws.Cells.Item(83, 9).Value = 2166.6667
$ws.Cells.Item(83, 10).Value = 3498
$ws.Cells.Item(83, 11).Value = 10833.3335
$ws.Cells.Item(83, 12).Value = 17490
$ws.Cells.Item(83, 13).Value = -5841.333500000001
$ws.Cells.Item(83, 14).Value = -27474
$ws.Cells.Item(107, 8).Value = 768.1875
$ws.Cells.Item(107, 9).Value = 582.7778
$ws.Cells.Item(107, 10).Value = 1006.5714
$ws.Cells.Item(107, 11).Value = 582.7778
$ws.Cells.Item(107, 12).Value = 1006.5714
$ws.Cells.Item(107, 13).Value = 1337.2222
$ws.Cells.Item(107, 14).Value = -4846.5714
$ws.Cells.Item(119, 8).Value = 29999
$ws.Cells.Item(119, 10).Value = 29999
$ws.Cells.Item(119, 12).Value = 29999
$ws.Cells.Item(119, 14).Value = -39675
$ws.Cells.Item(122, 8).Value = 3097.6775
$ws.Cells.Item(122, 9).Value = 2482.0454
$ws.Cells.Item(122, 11).Value = 7446.1362
$ws.Cells.Item(122, 13).Value = -4996.1362
$ws.Cells.Item(126, 8).Value = 16631.062
$ws.Cells.Item(126, 10).Value = 4199.6665
$ws.Cells.Item(126, 12).Value = 12598.9995
$ws.Cells.Item(126, 14).Value = -17538.9995
$ws.Cells.Item(132, 8).Value = 3357.4
$ws.Cells.Item(132, 10).Value = 3995.8572
$ws.Cells.Item(132, 12).Value = 11987.5716
$ws.Cells.Item(132, 14).Value = -17047.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 799.7273
$ws.Cells.Item(46, 9).Value = 730.875
$ws.Cells.Item(46, 10).Value = 983.3333
$ws.Cells.Item(46, 11).Value = 730.875
$ws.Cells.Item(46, 12).Value = 983.3333
$ws.Cells.Item(46, 13).Value = -542.875
$ws.Cells.Item(46, 14).Value = -1359.3333
$ws.Cells.Item(82, 8).Value = 8315.666999999999
$ws.Cells.Item(82, 9).Value = 9264
$ws.Cells.Item(82, 10).Value = 2151.5
$ws.Cells.Item(82, 11).Value = 9264
$ws.Cells.Item(82, 12).Value = 2151.5
$ws.Cells.Item(82, 13).Value = -8903
$ws.Cells.Item(82, 14).Value = -2873.5
$ws.Cells.Item(85, 8).Value = 8315.666999999999
$ws.Cells.Item(85, 9).Value = 9264
$ws.Cells.Item(85, 10).Value = 2151.5
$ws.Cells.Item(85, 11).Value = 9264
$ws.Cells.Item(85, 12).Value = 2151.5
$ws.Cells.Item(85, 13).Value = -8016
$ws.Cells.Item(85, 14).Value = -4647.5
$ws.Cells.Item(132, 8).Value = 3286.3635
$ws.Cells.Item(132, 9).Value = 3159.6956
$ws.Cells.Item(132, 10).Value = 3577.7
$ws.Cells.Item(132, 11).Value = 9479.086800000001
$ws.Cells.Item(132, 12).Value = 10733.1
$ws.Cells.Item(132, 13).Value = -6949.086800000001
$ws.Cells.Item(132, 14).Value = -15793.1
$ws.Cells.Item(136, 8).Value = 3633.5
$ws.Cells.Item(136, 9).Value = 3240.3
$ws.Cells.Item(136, 11).Value = 9720.900000000001
$ws.Cells.Item(136, 13).Value = -7170.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 6789.7856
$ws.Cells.Item(81, 9).Value = 9149.429
$ws.Cells.Item(81, 10).Value = 4430.143
$ws.Cells.Item(81, 11).Value = 18298.858
$ws.Cells.Item(81, 12).Value = 8860.286
$ws.Cells.Item(81, 13).Value = -17237.858
$ws.Cells.Item(81, 14).Value = -10982.286
$ws.Cells.Item(84, 8).Value = 6789.7856
$ws.Cells.Item(84, 9).Value = 9149.429
$ws.Cells.Item(84, 10).Value = 4430.143
$ws.Cells.Item(84, 11).Value = 91494.29000000001
$ws.Cells.Item(84, 12).Value = 44301.43
$ws.Cells.Item(84, 13).Value = -86190.29000000001
$ws.Cells.Item(84, 14).Value = -54909.43
$ws.Cells.Item(122, 8).Value = 1906.68
$ws.Cells.Item(122, 9).Value = 1876.6818
$ws.Cells.Item(122, 11).Value = 5630.0454
$ws.Cells.Item(122, 13).Value = -3180.0454
$ws.Cells.Item(126, 8).Value = 12230.2
$ws.Cells.Item(126, 9).Value = 12230.2
$ws.Cells.Item(126, 11).Value = 36690.60000000001
$ws.Cells.Item(126, 13).Value = -34220.60000000001
$ws.Cells.Item(132, 8).Value = 2373.0625
$ws.Cells.Item(132, 9).Value = 2274.48
$ws.Cells.Item(132, 11).Value = 6823.440000000001
$ws.Cells.Item(132, 13).Value = -4293.440000000001
$ws.Cells.Item(136, 8).Value = 2011.2307
$ws.Cells.Item(136, 9).Value = 1363.625
$ws.Cells.Item(136, 11).Value = 4090.875
$ws.Cells.Item(136, 13).Value = -1540.875
